$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (regenerated s_val data, filtered save games)
$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 16.44336272706294

# Row 3 (regenerated s_val data, filtered save games)
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 87981.0709163148
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 88229.71468683209
